$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")

# Update the "repaymentstrategy" value (B17): was "Mifos style", now the new
# periodic/upfront scenario label, carrying over the header-row formatting
# (same style as B1) instead of the old plain text style.
$ws.Range("B1").Copy()
$ws.Range("B17").PasteSpecial(-4122)
$ws.Range("B17").Value = "Penalties, Fees, Interest, Principal order"

# Move the active selection to B17 (this also clears the previous scrolled
# topLeftCell="A10" view position, restoring the default top-left scroll).
$ws.Range("B17").Select()
